# Refresh the cryptos snapshot sheet (Price / Volume(1h) columns) to
# match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose refreshed value is a plain decimal number (e.g.
# 0.999 / 300.25 / 2.20). Excel's smart-entry would silently coerce
# these into numeric literals (dropping the significant trailing
# zero / precision), so they are pinned to Text format first, same
# as how the original values were stored (general/inline text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Refresh Price (D) and Volume(1h) (E) columns with the new snapshot.
$ws.Range("D2").Value = '42.617.99'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '2.288.19'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '300.25'
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = '98.49'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("E7").Value = '  +0.85%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.508'
$ws.Range("E9").Value = '  +3.72%  '
$ws.Range("D10").Value = '35.53'
$ws.Range("E10").Value = '  +7.23%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("E12").Value = '  +2.11%  '
$ws.Range("E13").Value = '  +12.48%  '
$ws.Range("D14").Value = '6.78'
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("D15").Value = '2.648.20'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = '2.290.13'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '0.796'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '42.527.98'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("E19").Value = '  +5.79%  '
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("D21").Value = '0.0₃0895'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '67.64'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").Value = '234.87'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  +12.29%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = '24.42'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("D28").Value = '2.35'
$ws.Range("E28").Value = '  +13.47%  '
$ws.Range("D29").Value = '167.27'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '34.01'
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").Value = '9.10'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = '4.95'
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("E34").Value = '  -3.15%  '
$ws.Range("D35").Value = '17.13'
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("E36").Value = '  +3.45%  '
$ws.Range("D37").Value = '0.0687'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").Value = '2.81'
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '1.76'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").Value = '1.980.55'
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("E44").Value = '  -5.36%  '
$ws.Range("D45").Value = '10.05'
$ws.Range("E45").Value = '  +4.75%  '
$ws.Range("D46").Value = '17.38'
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").Value = '2.85'
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("D48").Value = '55.37'
$ws.Range("E48").Value = '  +6.25%  '
$ws.Range("D49").Value = '2.517.31'
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D50").Value = '1.51'
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("D51").Value = '4.48'
$ws.Range("E51").Value = '  +0.25%  '

# Drop back to the default (unstyled) cell style so the temporary
# Text number-format above doesn't leave a residual style diff.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
